# Clean codes for submission
# Replace hyphen-separated CI ranges with comma-separated ranges, and round p-values to 3 decimals.
$d = $word.ActiveDocument

$d.Content.Find.Execute("0.0124", $true, $false, $false, $false, $false, $true, 1, $false, "0.012", 2) | Out-Null
$d.Content.Find.Execute("1.18 (1.01-1.38)", $true, $false, $false, $false, $false, $true, 1, $false, "1.18 (1.01, 1.38)", 2) | Out-Null
$d.Content.Find.Execute("0.0340", $true, $false, $false, $false, $false, $true, 1, $false, "0.034", 2) | Out-Null
$d.Content.Find.Execute("1.21 (1.05-1.41)", $true, $false, $false, $false, $false, $true, 1, $false, "1.21 (1.05, 1.41)", 2) | Out-Null
$d.Content.Find.Execute("0.0106", $true, $false, $false, $false, $false, $true, 1, $false, "0.011", 2) | Out-Null
$d.Content.Find.Execute("1.16 (1.05-1.28)", $true, $false, $false, $false, $false, $true, 1, $false, "1.16 (1.05, 1.28)", 2) | Out-Null
$d.Content.Find.Execute("0.0040", $true, $false, $false, $false, $false, $true, 1, $false, "0.004", 2) | Out-Null
$d.Content.Find.Execute("1.12 (1.02-1.22)", $true, $false, $false, $false, $false, $true, 1, $false, "1.12 (1.02, 1.22)", 2) | Out-Null
$d.Content.Find.Execute("0.0130", $true, $false, $false, $false, $false, $true, 1, $false, "0.013", 2) | Out-Null
$d.Content.Find.Execute("0.0213", $true, $false, $false, $false, $false, $true, 1, $false, "0.021", 2) | Out-Null
$d.Content.Find.Execute("1.16 (1-1.36)", $true, $false, $false, $false, $false, $true, 1, $false, "1.16 (1, 1.36)", 2) | Out-Null
$d.Content.Find.Execute("0.0514", $true, $false, $false, $false, $false, $true, 1, $false, "0.051", 2) | Out-Null
$d.Content.Find.Execute("1.16 (0.99-1.34)", $true, $false, $false, $false, $false, $true, 1, $false, "1.16 (0.99, 1.34)", 2) | Out-Null
$d.Content.Find.Execute("0.0601", $true, $false, $false, $false, $false, $true, 1, $false, "0.060", 2) | Out-Null
$d.Content.Find.Execute("1.17 (1.06-1.29)", $true, $false, $false, $false, $false, $true, 1, $false, "1.17 (1.06, 1.29)", 2) | Out-Null
$d.Content.Find.Execute("0.0024", $true, $false, $false, $false, $false, $true, 1, $false, "0.002", 2) | Out-Null
$d.Content.Find.Execute("0.0144", $true, $false, $false, $false, $false, $true, 1, $false, "0.014", 2) | Out-Null
$d.Content.Find.Execute("0.0129", $true, $false, $false, $false, $false, $true, 1, $false, "0.013", 2) | Out-Null
$d.Content.Find.Execute("1.19 (1.03-1.39)", $true, $false, $false, $false, $false, $true, 1, $false, "1.19 (1.03, 1.39)", 2) | Out-Null
$d.Content.Find.Execute("0.0220", $true, $false, $false, $false, $false, $true, 1, $false, "0.022", 2) | Out-Null
$d.Content.Find.Execute("1.12 (0.96-1.3)", $true, $false, $false, $false, $false, $true, 1, $false, "1.12 (0.96, 1.3)", 2) | Out-Null
$d.Content.Find.Execute("0.1504", $true, $false, $false, $false, $false, $true, 1, $false, "0.150", 2) | Out-Null
$d.Content.Find.Execute("1.17 (1.06-1.3)", $true, $false, $false, $false, $false, $true, 1, $false, "1.17 (1.06, 1.3)", 2) | Out-Null
$d.Content.Find.Execute("0.0015", $true, $false, $false, $false, $false, $true, 1, $false, "0.001", 2) | Out-Null
$d.Content.Find.Execute("1.13 (1.03-1.23)", $true, $false, $false, $false, $false, $true, 1, $false, "1.13 (1.03, 1.23)", 2) | Out-Null
$d.Content.Find.Execute("0.0094", $true, $false, $false, $false, $false, $true, 1, $false, "0.009", 2) | Out-Null
$d.Content.Find.Execute("1.03 (0.99-1.06)", $true, $false, $false, $false, $false, $true, 1, $false, "1.03 (0.99, 1.06)", 2) | Out-Null
$d.Content.Find.Execute("0.1682", $true, $false, $false, $false, $false, $true, 1, $false, "0.168", 2) | Out-Null
$d.Content.Find.Execute("1.03 (0.99-1.07)", $true, $false, $false, $false, $false, $true, 1, $false, "1.03 (0.99, 1.07)", 2) | Out-Null
$d.Content.Find.Execute("0.1393", $true, $false, $false, $false, $false, $true, 1, $false, "0.139", 2) | Out-Null
$d.Content.Find.Execute("1.02 (0.99-1.06)", $true, $false, $false, $false, $false, $true, 1, $false, "1.02 (0.99, 1.06)", 2) | Out-Null
$d.Content.Find.Execute("0.2586", $true, $false, $false, $false, $false, $true, 1, $false, "0.259", 2) | Out-Null
$d.Content.Find.Execute("0.2085", $true, $false, $false, $false, $false, $true, 1, $false, "0.208", 2) | Out-Null
$d.Content.Find.Execute("1.2 (0.95-1.52)", $true, $false, $false, $false, $false, $true, 1, $false, "1.2 (0.95, 1.52)", 2) | Out-Null
$d.Content.Find.Execute("0.1148", $true, $false, $false, $false, $false, $true, 1, $false, "0.115", 2) | Out-Null
$d.Content.Find.Execute("1.28 (1.03-1.58)", $true, $false, $false, $false, $false, $true, 1, $false, "1.28 (1.03, 1.58)", 2) | Out-Null
$d.Content.Find.Execute("0.0273", $true, $false, $false, $false, $false, $true, 1, $false, "0.027", 2) | Out-Null
$d.Content.Find.Execute("1.1 (0.93-1.29)", $true, $false, $false, $false, $false, $true, 1, $false, "1.1 (0.93, 1.29)", 2) | Out-Null
$d.Content.Find.Execute("0.2590", $true, $false, $false, $false, $false, $true, 1, $false, "0.259", 2) | Out-Null
$d.Content.Find.Execute("1.08 (0.93-1.26)", $true, $false, $false, $false, $false, $true, 1, $false, "1.08 (0.93, 1.26)", 2) | Out-Null
$d.Content.Find.Execute("0.2955", $true, $false, $false, $false, $false, $true, 1, $false, "0.296", 2) | Out-Null
$d.Content.Find.Execute("0.2931", $true, $false, $false, $false, $false, $true, 1, $false, "0.293", 2) | Out-Null
$d.Content.Find.Execute("1.21 (0.96-1.53)", $true, $false, $false, $false, $false, $true, 1, $false, "1.21 (0.96, 1.53)", 2) | Out-Null
$d.Content.Find.Execute("0.0995", $true, $false, $false, $false, $false, $true, 1, $false, "0.100", 2) | Out-Null
$d.Content.Find.Execute("1.23 (0.99-1.53)", $true, $false, $false, $false, $false, $true, 1, $false, "1.23 (0.99, 1.53)", 2) | Out-Null
$d.Content.Find.Execute("0.0619", $true, $false, $false, $false, $false, $true, 1, $false, "0.062", 2) | Out-Null
$d.Content.Find.Execute("0.2566", $true, $false, $false, $false, $false, $true, 1, $false, "0.257", 2) | Out-Null
$d.Content.Find.Execute("1.09 (0.94-1.27)", $true, $false, $false, $false, $false, $true, 1, $false, "1.09 (0.94, 1.27)", 2) | Out-Null
$d.Content.Find.Execute("0.2570", $true, $false, $false, $false, $false, $true, 1, $false, "0.257", 2) | Out-Null
$d.Content.Find.Execute("0.1899", $true, $false, $false, $false, $false, $true, 1, $false, "0.190", 2) | Out-Null
$d.Content.Find.Execute("1.22 (0.97-1.54)", $true, $false, $false, $false, $false, $true, 1, $false, "1.22 (0.97, 1.54)", 2) | Out-Null
$d.Content.Find.Execute("0.0863", $true, $false, $false, $false, $false, $true, 1, $false, "0.086", 2) | Out-Null
$d.Content.Find.Execute("1.27 (1.02-1.58)", $true, $false, $false, $false, $false, $true, 1, $false, "1.27 (1.02, 1.58)", 2) | Out-Null
$d.Content.Find.Execute("0.0307", $true, $false, $false, $false, $false, $true, 1, $false, "0.031", 2) | Out-Null
$d.Content.Find.Execute("1.07 (0.91-1.26)", $true, $false, $false, $false, $false, $true, 1, $false, "1.07 (0.91, 1.26)", 2) | Out-Null
$d.Content.Find.Execute("0.4006", $true, $false, $false, $false, $false, $true, 1, $false, "0.401", 2) | Out-Null
$d.Content.Find.Execute("1.1 (0.95-1.28)", $true, $false, $false, $false, $false, $true, 1, $false, "1.1 (0.95, 1.28)", 2) | Out-Null
$d.Content.Find.Execute("0.2073", $true, $false, $false, $false, $false, $true, 1, $false, "0.207", 2) | Out-Null
$d.Content.Find.Execute("1.04 (0.98-1.09)", $true, $false, $false, $false, $false, $true, 1, $false, "1.04 (0.98, 1.09)", 2) | Out-Null
$d.Content.Find.Execute("0.1966", $true, $false, $false, $false, $false, $true, 1, $false, "0.197", 2) | Out-Null
$d.Content.Find.Execute("1.05 (0.99-1.1)", $true, $false, $false, $false, $false, $true, 1, $false, "1.05 (0.99, 1.1)", 2) | Out-Null
$d.Content.Find.Execute("0.0952", $true, $false, $false, $false, $false, $true, 1, $false, "0.095", 2) | Out-Null
$d.Content.Find.Execute("1.05 (0.99-1.11)", $true, $false, $false, $false, $false, $true, 1, $false, "1.05 (0.99, 1.11)", 2) | Out-Null
$d.Content.Find.Execute("0.0783", $true, $false, $false, $false, $false, $true, 1, $false, "0.078", 2) | Out-Null
